$d = $word.ActiveDocument

# Locate the "LOQ4056 ... (Requisito fraco)" paragraph; the three
# paragraphs right after it (a blank paragraph, a blank page-break
# paragraph, and the "(c) 2020 ... Attribution" footer notice) are the
# ones being dropped from the end of the document.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains("Requisito fraco")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Remove back-to-front so earlier indices stay valid while deleting.
    $d.Paragraphs.Item($targetIndex + 3).Range.Delete()
    $d.Paragraphs.Item($targetIndex + 2).Range.Delete()
    $d.Paragraphs.Item($targetIndex + 1).Range.Delete()
}
